$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 - SALDO AWAL (only the Saldo value changes)
$ws.Cells.Item(2, 5).Value = 76939992.8

# Row 3
Set-TextCell 3 1 "04/06/2024"
Set-TextCell 3 2 "PAY KARTU KREDIT 5498460025096193"
$ws.Cells.Item(3, 3).Value = 33224480
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 43715512.8

# Row 4
Set-TextCell 4 1 "10/06/2024"
Set-TextCell 4 2 "TRF DARI NI KETUT KARSINI BANK"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 2233833
$ws.Cells.Item(4, 5).Value = 45949345.8

# Row 5
Set-TextCell 5 1 "10/06/2024"
Set-TextCell 5 2 "TRF DARI NI KETUT KARSINI BANK"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 8250000
$ws.Cells.Item(5, 5).Value = 54199345.8

# Row 6
Set-TextCell 6 1 "18/06/2024"
Set-TextCell 6 2 "QR PAYMENT 16.09.48 MONSIEUR"
$ws.Cells.Item(6, 3).Value = 130900
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 54068445.8

# Row 7 (previously SALDO AKHIR, now a transaction row)
Set-TextCell 7 1 "19/06/2024"
Set-TextCell 7 2 "PB KE NI PUTU SAWITRI 9996413192 08"
$ws.Cells.Item(7, 3).Value = 8025000
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 46043445.8

# Row 8 (new)
Set-TextCell 8 1 "20/06/2024"
Set-TextCell 8 2 "QR PAYMENT 16.22.27 SOULSHINE"
$ws.Cells.Item(8, 3).Value = 1881550
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 44161895.8

# Row 9 (new)
Set-TextCell 9 1 "30/06/2024"
Set-TextCell 9 2 "PENDAPATAN BUNGA"
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 10350
$ws.Cells.Item(9, 5).Value = 44172245.8

# Row 10 (new)
Set-TextCell 10 1 "30/06/2024"
Set-TextCell 10 2 "PAJAK ATAS BUNGA"
$ws.Cells.Item(10, 3).Value = 2070
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 44170175.8

# Row 11 (new) - SALDO AKHIR, moved from row 7
Set-TextCell 11 2 "SALDO AKHIR"
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 44170175.8
